# Update countries & provincias Spain
# Applies the 20-Abril-2020 15:52 data refresh to the "Pais" sheet:
#   - Bielorrusia's case count overtakes Catar/Ucrania/Malasia, so it moves
#     up a row (A44..A47 relabelled, B..H values shifted accordingly)
#   - Numeric case/death/recovery counters refreshed for several countries
#   - Footer timestamp bumped from 15:22 to 15:52

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 15:52"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Range("B4").Value = 765613
$ws.Range("C4").Value = 977
$ws.Range("D4").Value = 71253
$ws.Range("E4").Value = 653740
$ws.Range("F4").Value = 13566
$ws.Range("G4").Value = 45
$ws.Range("H4").Value = 40620

# --- Row 15: Brasil -----------------------------------------------------
$ws.Range("F15").Value = 7919

# --- Rows 44-47: Bielorrusia overtakes Catar/Ucrania/Malasia ----------
# Row 44 becomes Bielorrusia with fresh data
$ws.Range("A44").Value = "Bielorrusia"
$ws.Range("B44").Value = 6264
$ws.Range("C44").Value = 1485
$ws.Range("D44").Value = 514
$ws.Range("E44").Value = 5699
$ws.Range("F44").Value = 92
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 51

# Row 45 becomes Catar (previously row 44's data, unchanged)
$ws.Range("A45").Value = "Catar"
$ws.Range("B45").Value = 6015
$ws.Range("C45").Value = 567
$ws.Range("D45").Value = 555
$ws.Range("E45").Value = 5451
$ws.Range("F45").Value = 37
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 9

# Row 46 becomes Ucrania (previously row 45's data, unchanged)
$ws.Range("A46").Value = "Ucrania"
$ws.Range("B46").Value = 5710
$ws.Range("C46").Value = 261
$ws.Range("D46").Value = 359
$ws.Range("E46").Value = 5200
$ws.Range("F46").Value = 45
$ws.Range("G46").Value = 10
$ws.Range("H46").Value = 151

# Row 47 becomes Malasia (previously row 46's data, unchanged)
$ws.Range("A47").Value = "Malasia"
$ws.Range("B47").Value = 5425
$ws.Range("C47").Value = 36
$ws.Range("D47").Value = 3295
$ws.Range("E47").Value = 2041
$ws.Range("F47").Value = 45
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 89

# --- Row 67: Islandia -----------------------------------------------
$ws.Range("B67").Value = 1773
$ws.Range("C67").Value = 2
$ws.Range("D67").Value = 1362
$ws.Range("E67").Value = 402
$ws.Range("F67").Value = 4

# --- Row 93: Libano ----------------------------------------------------
$ws.Range("D93").Value = 103
$ws.Range("E93").Value = 553

# --- Row 104: San Marino ------------------------------------------------
$ws.Range("B104").Value = 462
$ws.Range("C104").Value = 1
$ws.Range("D104").Value = 61

# --- Row 115: Sri Lanka --------------------------------------------------
$ws.Range("B115").Value = 304
$ws.Range("C115").Value = 33
$ws.Range("D115").Value = 98

# --- Row 159: Uganda -----------------------------------------------------
$ws.Range("D159").Value = 38
$ws.Range("E159").Value = 17
